$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply AutoFilter to the data range:
#  - Column A (ObsID): only rows starting with "2023"
#  - Column D (CM3): only rows with value between 35 and 55 (inclusive)
$rng = $ws.Range("A1:O78")
$rng.AutoFilter(1, "2023*")
$rng.AutoFilter(4, ">=35", 1, "<=55")

# Excel records the filtered range as a hidden workbook-level defined name
# scoped to this sheet.
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "='Catalog-static'!`$A`$1:`$O`$78")
$fdb.Visible = $false

# Update the current selection/view on the sheet.
$ws.Range("J80").Select()
